$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 125001940
$ws.Range("I55").Value = 200002530
$ws.Range("K55").Value = 200002530
$ws.Range("M55").Value = -200002316
$ws.Range("H112").Value = 5424.4116
$ws.Range("J112").Value = 5913.871
$ws.Range("L112").Value = 17741.613
$ws.Range("N112").Value = -19957.613
$ws.Range("H113").Value = 2033.05
$ws.Range("I113").Value = 1653.8889
$ws.Range("J113").Value = 2343.2727
$ws.Range("K113").Value = 1653.8889
$ws.Range("L113").Value = 2343.2727
$ws.Range("M113").Value = 1600.1111
$ws.Range("N113").Value = -8851.2727
$ws.Range("H135").Value = 966.8182
$ws.Range("I135").Value = 513.0645
$ws.Range("J135").Value = 8000
$ws.Range("K135").Value = 4617.5805
$ws.Range("L135").Value = 72000
$ws.Range("M135").Value = -2082.5805
$ws.Range("N135").Value = -77070
$ws.Range("H138").Value = 1927.2046
$ws.Range("I138").Value = 1452.1945
$ws.Range("J138").Value = 2256.0576
$ws.Range("K138").Value = 4356.583500000001
$ws.Range("L138").Value = 6768.1728
$ws.Range("M138").Value = 783.4164999999994
$ws.Range("N138").Value = -17048.1728

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2949.6
$ws.Range("I45").Value = 2862
$ws.Range("J45").Value = 3300
$ws.Range("K45").Value = 2862
$ws.Range("L45").Value = 3300
$ws.Range("M45").Value = -2485
$ws.Range("N45").Value = -4054
$ws.Range("H132").Value = 2936.0952
$ws.Range("I132").Value = 2343.7754
$ws.Range("J132").Value = 5009.2144
$ws.Range("K132").Value = 7031.3262
$ws.Range("L132").Value = 15027.6432
$ws.Range("M132").Value = -4501.3262
$ws.Range("N132").Value = -20087.6432

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4612.174
$ws.Range("I134").Value = 4240.3076
$ws.Range("J134").Value = 5095.6
$ws.Range("K134").Value = 12720.9228
$ws.Range("L134").Value = 15286.8
$ws.Range("M134").Value = -10185.9228
$ws.Range("N134").Value = -20356.8

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3365.4744
$ws.Range("I31").Value = 1082.1428
$ws.Range("K31").Value = 1082.1428
$ws.Range("M31").Value = -787.1428000000001
$ws.Range("H34").Value = 3365.4744
$ws.Range("I34").Value = 1082.1428
$ws.Range("K34").Value = 1082.1428
$ws.Range("M34").Value = -880.1428000000001
$ws.Range("H58").Value = 1101.9546
$ws.Range("I58").Value = 789.48
$ws.Range("J58").Value = 1513.1052
$ws.Range("K58").Value = 789.48
$ws.Range("L58").Value = 1513.1052
$ws.Range("M58").Value = -586.48
$ws.Range("N58").Value = -1919.1052
$ws.Range("H99").Value = 2000
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H135").Value = 40660
$ws.Range("J135").Value = 40660
$ws.Range("L135").Value = 40660
$ws.Range("N135").Value = -50800
$ws.Range("H136").Value = 1101.9546
$ws.Range("I136").Value = 789.48
$ws.Range("J136").Value = 1513.1052
$ws.Range("K136").Value = 2368.44
$ws.Range("L136").Value = 4539.3156
$ws.Range("M136").Value = 181.5599999999999
$ws.Range("N136").Value = -9639.3156

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 7778.4165
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 7778.4165
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 23335.2495
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -24707.2495
$ws.Range("H65").Value = 7778.4165
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 7778.4165
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 70005.7485
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -76869.7485
$ws.Range("H122").Value = 4173.069
$ws.Range("I122").Value = 478.6111
$ws.Range("J122").Value = 10218.546
$ws.Range("K122").Value = 4307.4999
$ws.Range("L122").Value = 91966.914
$ws.Range("M122").Value = -1857.4999
$ws.Range("N122").Value = -96866.914
$ws.Range("H131").Value = 2852.0938
$ws.Range("J131").Value = 3240.7454
$ws.Range("L131").Value = 9722.236199999999
$ws.Range("N131").Value = -19802.2362
$ws.Range("H132").Value = 2709.3684
$ws.Range("I132").Value = 2828.2222
$ws.Range("J132").Value = 2672.4827
$ws.Range("K132").Value = 25453.9998
$ws.Range("L132").Value = 24052.3443
$ws.Range("M132").Value = -22923.9998
$ws.Range("N132").Value = -29112.3443

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 59.666668
$ws.Range("I2").Value = 69.5
$ws.Range("J2").Value = 40
$ws.Range("K2").Value = 69.5
$ws.Range("L2").Value = 40
$ws.Range("M2").Value = 43.5
$ws.Range("N2").Value = -266
$ws.Range("H29").Value = 5700.8887
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 5700.8887
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 5700.8887
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -6280.8887
$ws.Range("H70").Value = 8579.482
$ws.Range("I70").Value = 9409.409
$ws.Range("K70").Value = 9409.409
$ws.Range("M70").Value = -9139.409
$ws.Range("H73").Value = 8579.482
$ws.Range("I73").Value = 9409.409
$ws.Range("K73").Value = 9409.409
$ws.Range("M73").Value = -8473.409
$ws.Range("H113").Value = 127374.875
$ws.Range("I113").Value = 168833.17
$ws.Range("K113").Value = 168833.17
$ws.Range("M113").Value = -166663.17
$ws.Range("H126").Value = 1966.8334
$ws.Range("I126").Value = 1967
$ws.Range("J126").Value = 1966.6666
$ws.Range("K126").Value = 5901
$ws.Range("L126").Value = 5899.9998
$ws.Range("M126").Value = -3431
$ws.Range("N126").Value = -10839.9998
$ws.Range("H132").Value = 4382.1177
$ws.Range("I132").Value = 3677.2856
$ws.Range("K132").Value = 11031.8568
$ws.Range("M132").Value = -8501.856800000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6230.25
$ws.Range("I7").Value = 6668
$ws.Range("K7").Value = 6668
$ws.Range("M7").Value = -6556
$ws.Range("H40").Value = 126825.25
$ws.Range("I40").Value = 168100.33
$ws.Range("K40").Value = 168100.33
$ws.Range("M40").Value = -167964.33
$ws.Range("H122").Value = 3630.6667
$ws.Range("I122").Value = 3164.8462
$ws.Range("J122").Value = 4063.2144
$ws.Range("K122").Value = 9494.5386
$ws.Range("L122").Value = 12189.6432
$ws.Range("M122").Value = -7044.5386
$ws.Range("N122").Value = -17089.6432
$ws.Range("H126").Value = 6230.25
$ws.Range("I126").Value = 6668
$ws.Range("K126").Value = 20004
$ws.Range("M126").Value = -17534
$ws.Range("H132").Value = 2843.3262
$ws.Range("I132").Value = 3099.3845
$ws.Range("J132").Value = 2510.45
$ws.Range("K132").Value = 9298.1535
$ws.Range("L132").Value = 7531.349999999999
$ws.Range("M132").Value = -6768.1535
$ws.Range("N132").Value = -12591.35
$ws.Range("H136").Value = 5210189
$ws.Range("I136").Value = 2027.25
$ws.Range("J136").Value = 20834676
$ws.Range("K136").Value = 6081.75
$ws.Range("L136").Value = 62504028
$ws.Range("M136").Value = -3531.75
$ws.Range("N136").Value = -62509128

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 60001.5
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 60001.5
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 60001.5
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -60225.5
$ws.Range("H18").Value = 70007
$ws.Range("J18").Value = 70007
$ws.Range("L18").Value = 70007
$ws.Range("N18").Value = -70353
$ws.Range("H122").Value = 2387.5454
$ws.Range("I122").Value = 2380.7693
$ws.Range("J122").Value = 2412.7144
$ws.Range("K122").Value = 7142.3079
$ws.Range("L122").Value = 7238.1432
$ws.Range("M122").Value = -4692.3079
$ws.Range("N122").Value = -12138.1432
$ws.Range("H126").Value = 829.3889
$ws.Range("I126").Value = 784.1
$ws.Range("J126").Value = 886
$ws.Range("K126").Value = 2352.3
$ws.Range("L126").Value = 2658
$ws.Range("M126").Value = 117.6999999999998
$ws.Range("N126").Value = -7598
$ws.Range("H132").Value = 4239
$ws.Range("I132").Value = 5302.7144
$ws.Range("J132").Value = 2377.5
$ws.Range("K132").Value = 15908.1432
$ws.Range("L132").Value = 7132.5
$ws.Range("M132").Value = -13378.1432
$ws.Range("N132").Value = -12192.5
$ws.Range("H136").Value = 2984.625
$ws.Range("I136").Value = 2697.8857
$ws.Range("J136").Value = 3756.6155
$ws.Range("K136").Value = 8093.657099999999
$ws.Range("L136").Value = 11269.8465
$ws.Range("M136").Value = -5543.657099999999
$ws.Range("N136").Value = -16369.8465

Write-Output "Applied all updates"